$d = $word.ActiveDocument

# Update the date/title line (first paragraph in the document)
$d.Paragraphs.Item(1).Range.Text = "2025-06-13 Friday"

$tbl = $d.Tables.Item(1)

# Row 1
# 35÷9= -> 91÷8=
$tbl.Cell(1, 1).Range.Text = "91÷8="
# 19÷7= -> 60÷8=
$tbl.Cell(1, 2).Range.Text = "60÷8="
# 94÷5= -> 46÷5=
$tbl.Cell(1, 3).Range.Text = "46÷5="
# 45÷6= -> 18÷2=
$tbl.Cell(1, 4).Range.Text = "18÷2="
# 29÷6= -> 87÷5=
$tbl.Cell(1, 5).Range.Text = "87÷5="

# Row 5
# 63÷5= -> 72÷2=
$tbl.Cell(5, 1).Range.Text = "72÷2="
# 67÷6= -> 55÷7=
$tbl.Cell(5, 2).Range.Text = "55÷7="
# 48÷2= -> 59÷5=
$tbl.Cell(5, 3).Range.Text = "59÷5="
# 25÷8= -> 67÷2=
$tbl.Cell(5, 4).Range.Text = "67÷2="
# 54÷3= -> 10÷4=
$tbl.Cell(5, 5).Range.Text = "10÷4="

# Row 9
# 88÷6= -> 63÷5=
$tbl.Cell(9, 1).Range.Text = "63÷5="
# 95÷2= -> 75÷5=
$tbl.Cell(9, 2).Range.Text = "75÷5="
# 33÷6= -> 28÷2=
$tbl.Cell(9, 3).Range.Text = "28÷2="
# 18÷8= -> 88÷9=
$tbl.Cell(9, 4).Range.Text = "88÷9="
# 88÷3= -> 86÷6=
$tbl.Cell(9, 5).Range.Text = "86÷6="

# Row 13
# 31÷3= -> 21÷4=
$tbl.Cell(13, 1).Range.Text = "21÷4="
# 32÷2= -> 31÷9=
$tbl.Cell(13, 2).Range.Text = "31÷9="
# 90÷5= -> 53÷3=
$tbl.Cell(13, 3).Range.Text = "53÷3="
# 96÷4= -> 15÷7=
$tbl.Cell(13, 4).Range.Text = "15÷7="
# 89÷2= -> 79÷4=
$tbl.Cell(13, 5).Range.Text = "79÷4="

# Row 17
# 50÷9= -> 45÷8=
$tbl.Cell(17, 1).Range.Text = "45÷8="
# 24÷7= -> 85÷4=
$tbl.Cell(17, 2).Range.Text = "85÷4="
# 27÷9= -> 22÷9=
$tbl.Cell(17, 3).Range.Text = "22÷9="
# 63÷3= -> 34÷6=
$tbl.Cell(17, 4).Range.Text = "34÷6="
# 67÷2= -> 35÷3=
$tbl.Cell(17, 5).Range.Text = "35÷3="

